$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "popularity" column
$ws.Range("S1").Value = "popularity"

# Copy header style (bold, centered, bordered) from R1 onto S1
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate popularity values for data rows 2-51
$popularity = @{
    2 = 95
    3 = 94
    4 = 86
    5 = 91
    6 = 91
    7 = 91
    8 = 91
    9 = 93
    10 = 91
    11 = 96
    12 = 92
    13 = 89
    14 = 95
    15 = 88
    16 = 89
    17 = 86
    18 = 92
    19 = 90
    20 = 92
    21 = 89
    22 = 89
    23 = 86
    24 = 86
    25 = 87
    26 = 84
    27 = 85
    28 = 100
    29 = 84
    30 = 88
    31 = 88
    32 = 84
    33 = 91
    34 = 87
    35 = 84
    36 = 86
    37 = 84
    38 = 88
    39 = 78
    40 = 87
    41 = 85
    42 = 88
    43 = 87
    44 = 89
    45 = 83
    46 = 88
    47 = 86
    48 = 84
    49 = 90
    50 = 82
    51 = 78
}

foreach ($row in $popularity.Keys) {
    $ws.Cells.Item($row, 19).Value = $popularity[$row]
}

